$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D5").Value = "연속 신호의 샘플링"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2022/01/14/sampling_CT_to_DT.html"

$ws.Range("D36").Value = "Time Series Representation Learning"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/348"

$ws.Range("D52").Value = "html_table()과 invalid multibyte string"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2384&utm_source=rss&utm_medium=rss&utm_campaign=html_table%25ea%25b3%25bc-invalid-multibyte-string"
